$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("Bradley Beal", "PG,SG,SF", "Phoenix Suns"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("Zion Williamson", "SF,PF,C", "New Orleans Pelicans"),
    @("Quentin Grimes", "SG,SF", "Philadelphia 76ers"),
    @("RJ Barrett", "SG,SF,PF", "Toronto Raptors"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Jabari Smith Jr.", "PF,C", "Houston Rockets"),
    @("Draymond Green", "PF,C", "Golden State Warriors"),
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Donte DiVincenzo", "PG,SG,SF", "Minnesota Timberwolves"),
    @("Amen Thompson", "PG,SG,SF,PF", "Houston Rockets"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
